# Commit message: "Create CSV file with Classic answers"
# The edit corrects the punctuation of one of the Magic 8 Ball answers on the
# "Magic Answers" sheet: "Yes – definitely" (en dash) becomes
# "Yes, definitely" (comma).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Magic Answers")

# B20 holds the "Yes – definitely" answer text; fix the punctuation.
$ws.Range("B20").Value = "Yes, definitely"

# Reflect the final selection recorded in the saved file.
$ws.Activate()
$ws.Range("B22").Select()
